$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(129, 1).Value = "Sr.Software Developer(Golang, SQL, Postgres)"
$ws.Cells.Item(129, 2).Value = "https://www.dice.com/job-detail/860a5fe9-3cd6-43de-a64a-b4435aa147e5"
$ws.Cells.Item(129, 3).Value = "Remote"
$ws.Cells.Item(129, 4).Value = "Contract"
$ws.Cells.Item(129, 5).Value = "Depends on Experience"
$ws.Cells.Item(129, 6).Value = "STAND 8"

$ws.Cells.Item(130, 1).Value = "Golang Developer"
$ws.Cells.Item(130, 2).Value = "https://www.dice.com/job-detail/ad009a2b-e4d0-4713-a240-7edb44f0a065"
$ws.Cells.Item(130, 3).Value = "McLean, Virginia"
$ws.Cells.Item(130, 4).Value = "Contract"
$ws.Cells.Item(130, 5).Value = "`$50 - `$60"
$ws.Cells.Item(130, 6).Value = "ICS Global Soft, Inc."
